$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in attendance ("P") for the three extra weeks (columns J, K, L) ---
$rng = $ws.Range("J2:L9")
$rng.Value = "P"
$rng.Font.Color = 5287936   # RGB(0, 176, 80) -> same green font used for "Present"
$rng.Font.Name = "Arial"
$rng.Font.Size = 12

# --- Move the active selection, as left behind by the editor ---
$ws.Range("K14").Select()

# --- Re-apply the Present/Absent conditional formatting (matches editing history
#     where the rule was re-created a few times before settling) ---
$cfRange = $ws.Range("B2:L9")
$cfRange.FormatConditions.Delete()

for ($i = 0; $i -lt 2; $i++) {
    $green = $cfRange.FormatConditions.Add(2, 3, 'B2="O"')
    $green.Interior.Color = 5296274   # FF92D050
    $red = $cfRange.FormatConditions.Add(2, 3, 'B2="X"')
    $red.Interior.Color = 255         # FFFF0000
    $cfRange.FormatConditions.Delete()
}

$greenFinal = $cfRange.FormatConditions.Add(2, 3, 'B2="O"')
$greenFinal.Interior.Color = 5296274
$redFinal = $cfRange.FormatConditions.Add(2, 3, 'B2="X"')
$redFinal.Interior.Color = 255
$redFinal.Priority = 1
$greenFinal.Priority = 2
